$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently holds the shared text string "R40". The edit changes it
# to hold the literal TEXT "1" (not the number 1), while leaving the cell's
# existing style/number-format untouched. A plain Value = "1" assignment
# would be auto-coerced to the number 1, so instead we stage the text "1" as
# a formula result in a scratch cell, copy it, and paste-special only the
# value into B11 - this records "1" as text without perturbing B11's style.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false
